{"js": "// Office.js (Word JS API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the change (per the diff):\n//  - Paragraph 1 (title block): date 12.02.25 -> 10.02.25, and the paper\n//    title line (after the manual line break) is replaced with the new title.\n//  - Paragraphs 2-17 (the body of the review) get their text replaced in\n//    place, one-for-one, with the new \"KANs\" review content.\n//  - Paragraphs 18-28 (old content about State Collapse / SC strategies /\n//    summary) are removed entirely (11 paragraphs deleted).\n//  - The final paragraph (the arxiv link) keeps its position but its text\n//    changes to the new URL.\n\nconst middleTexts = [\n  \"\u05de\u05d1\u05d5\u05d0:\",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05d0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05de\u05e6\u05d9\u05d2 \u05d7\u05e7\u05d9\u05e8\u05d4 \u05de\u05e2\u05de\u05d9\u05e7\u05d4 \u05e9\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3 (KANs), \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d7\u05d3\u05e9\u05e0\u05d9\u05ea \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc \u05de\u05e9\u05e4\u05d8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e9\u05d5\u05d5\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05de\u05d3\u05d5\u05e7\u05d3\u05e7 \u05d1\u05d9\u05df KANs \u05dc\u05d1\u05d9\u05df \u05e8\u05e9\u05ea\u05d5\u05ea MLPs \u05de\u05e1\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea, \u05d4\u05df \u05de\u05d1\u05d7\u05d9\u05e0\u05d4 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d5\u05d4\u05df \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea, \u05ea\u05d5\u05da \u05d4\u05ea\u05de\u05e7\u05d3\u05d5\u05ea \u05d1\u05d4\u05d9\u05d1\u05d8\u05d9\u05dd \u05db\u05de\u05d5 \u05d0\u05e7\u05e1\u05e4\u05e8\u05e1\u05d9\u05d1\u05e0\u05e1, \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d5\u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05ea \u05d0\u05d9\u05de\u05d5\u05df. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d1\u05e1\u05e1 \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05e8\u05db\u05d6\u05d9\u05d5\u05ea \u05d5\u05de\u05d0\u05de\u05ea \u05d0\u05d5\u05ea\u05df \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9\u05dd, \u05d5\u05d1\u05db\u05da \u05de\u05d4\u05d5\u05d5\u05d4 \u05ea\u05e8\u05d5\u05de\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05dc\u05ea\u05db\u05e0\u05d5\u05df \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1 \u05e9\u05d5\u05e0\u05d5\u05ea.\",\n  \"\u05d0\u05e7\u05e1\u05e4\u05e8\u05e1\u05d9\u05d1\u05e0\u05e1:\",\n  \"\u05d4\u05d9\u05e9\u05d2 \u05de\u05e8\u05db\u05d6\u05d9 \u05e9\u05dc \u05e2\u05d1\u05d5\u05d3\u05d4 \u05d6\u05d5 \u05d4\u05d5\u05d0 \u05d4\u05d4\u05d5\u05db\u05d7\u05d4 \u05d4\u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e9- KANs \u05d4\u05df \u05d1\u05e2\u05dc\u05d5\u05ea \u05d0\u05e7\u05e1\u05e4\u05e8\u05e1\u05d9\u05d1\u05e0\u05e1 \u05dc\u05e4\u05d7\u05d5\u05ea \u05db\u05de\u05d5 MLPs. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05db\u05dc MLP \u05de\u05d1\u05d5\u05e1\u05e1\u05ea ReLU \u05e0\u05d9\u05ea\u05df \u05f4\u05dc\u05de\u05e4\u05d5\u05ea\u05f4 \u05dc\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea KAN \u05de\u05e7\u05d1\u05d9\u05dc\u05d4, \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05d4 \u05e2\u05dc \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d5\u05dc\u05dc\u05d0 \u05d4\u05d2\u05d3\u05dc\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05d1\u05d2\u05d5\u05d3\u05dc \u05d4\u05e8\u05e9\u05ea. \u05de\u05e0\u05d2\u05d3, \u05d1\u05e2\u05d5\u05d3 \u05e9-KANs \u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d2\u05dd \u05e2\u05dc \u05d9\u05d3\u05d9 MLPs, \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e6\u05d9\u05d4 \u05d6\u05d5 \u05db\u05e8\u05d5\u05db\u05d4 \u05d1\u05e2\u05dc\u05d5\u05ea \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea: \u05de\u05e1\u05e4\u05e8 \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d2\u05d3\u05dc \u05e2\u05dd \u05d2\u05d5\u05d3\u05dc \u05d2\u05e8\u05d9\u05d3 (\u05de\u05e1\u05e4\u05e8 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05e2\u05d5\u05d2\u05df \u05d1\u05e1\u05e4\u05dc\u05d9\u05d9\u05df) \u05e9\u05dc \u05d4-KAN. \u05de\u05de\u05e6\u05d0 \u05d6\u05d4 \u05de\u05e8\u05de\u05d6 \u05e9-KANs \u05e2\u05e9\u05d5\u05d9\u05d5\u05ea \u05dc\u05d4\u05e6\u05d9\u05e2 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d9\u05e2\u05d9\u05dc\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05e2\u05d1\u05d5\u05e8 \u05e1\u05d5\u05d2\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea, \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05db\u05d0\u05e9\u05e8 \u05e0\u05e2\u05e9\u05d4 \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05de\u05d1\u05e0\u05d9 \u05d2\u05e8\u05d9\u05d3 \u05e2\u05d3\u05d9\u05e0\u05d9\u05dd.\",\n  \"\u05d4\u05de\u05d7\u05e7\u05e8 \u05de\u05e0\u05e6\u05dc \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05e7\u05d9\u05d9\u05de\u05d5\u05ea \u05e2\u05d1\u05d5\u05e8 MLPs \u05db\u05d3\u05d9 \u05dc\u05e7\u05d1\u05d5\u05e2 \u05e7\u05e6\u05d1\u05d9 \u05e7\u05d9\u05e8\u05d5\u05d1 \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05d1\u05d5\u05e8 KANs \u05d1\u05de\u05e8\u05d7\u05d1\u05d9\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05d0\u05d5\u05e0\u05dc\u05d9\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05db\u05de\u05d5 \u05de\u05e8\u05d7\u05d1 \u05e1\u05d5\u05d1\u05d5\u05dc\u05d1. \u05d4\u05d5\u05d0 \u05de\u05d3\u05d2\u05d9\u05dd \u05e9-KANs \u05de\u05e9\u05d9\u05d2\u05d5\u05ea \u05e7\u05e6\u05d1\u05d9 \u05e7\u05d9\u05e8\u05d5\u05d1 \u05d3\u05d5\u05de\u05d9\u05dd \u05d0\u05d5 \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d0\u05e9\u05e8 MLPs \u05d1\u05e9\u05e2\u05e8\u05d5\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea, \u05de\u05d4 \u05e9\u05de\u05d7\u05d6\u05e7 \u05d0\u05ea \u05d7\u05d5\u05e1\u05e0\u05df \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9.\",\n  \"\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05d8\u05d9\u05d9\u05ea \u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea (spectral bias):\",\n  \"\u05d0\u05d7\u05d3 \u05d4\u05d4\u05d1\u05d3\u05dc\u05d9\u05dd \u05d4\u05de\u05e8\u05db\u05d6\u05d9\u05d9\u05dd \u05d1\u05d9\u05df KANs \u05dc-MLPs \u05d4\u05de\u05d5\u05d3\u05d2\u05e9\u05d9\u05dd \u05d1\u05de\u05d0\u05de\u05e8 \u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d4\u05d4\u05d1\u05d3\u05dc \u05d1\u05d4\u05d8\u05d9\u05d4 \u05d4\u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea \u05e9\u05dc\u05d4\u05dd - \u05ea\u05d5\u05e4\u05e2\u05d4 \u05e9\u05d1\u05d4 \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05e0\u05d5\u05d8\u05d5\u05ea \u05dc\u05dc\u05de\u05d5\u05d3 \u05ea\u05d7\u05d9\u05dc\u05d4 \u05d1\u05ea\u05d3\u05e8\u05d9\u05dd \u05e0\u05de\u05d5\u05db\u05d9\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05d2\u05d9\u05dd \u05e0\u05d9\u05ea\u05d5\u05d7 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05d5\u05d0\u05de\u05e4\u05d9\u05e8\u05d9 \u05de\u05e4\u05d5\u05e8\u05d8, \u05d4\u05de\u05e8\u05d0\u05d4 \u05e9- KANs \u05e1\u05d5\u05d1\u05dc\u05d5\u05ea \u05e4\u05d7\u05d5\u05ea \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05de\u05d4\u05d8\u05d9\u05d4 \u05d6\u05d5.\",\n  \"\u05d4\u05d1\u05d3\u05dc \u05d6\u05d4 \u05de\u05d9\u05d5\u05d7\u05e1 \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05d4-B-spline \u05d5\u05dc\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d4\u05e7\u05d5\u05de\u05e4\u05d5\u05d6\u05d9\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9\u05ea \u05e9\u05dc KANs, \u05d4\u05de\u05d0\u05e4\u05e9\u05e8\u05d5\u05ea \u05dc\u05d4\u05df \u05dc\u05dc\u05de\u05d5\u05d3 \u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4 \u05d1\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05e8\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05e6\u05d9\u05e2\u05d5\u05ea \u05e9\u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05ea \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc KANs \u05e8\u05d3\u05d5\u05d3\u05d5\u05ea \u05d0\u05d7\u05d9\u05d3\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d9\u05d7\u05e1 \u05dc\u05ea\u05d3\u05e8\u05d9\u05dd \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05d1\u05d4\u05e9\u05d5\u05d5\u05d0\u05d4 \u05dc-MLPs, \u05e9\u05d1\u05d4\u05df \u05e0\u05e6\u05e4\u05d9\u05ea \u05d4\u05ea\u05db\u05e0\u05e1\u05d5\u05ea \u05de\u05d4\u05d9\u05e8\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05ea\u05d3\u05e8\u05d9\u05dd \u05e0\u05de\u05d5\u05db\u05d9\u05dd. \u05d4\u05d4\u05d8\u05d9\u05d4 \u05d4\u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea \u05d4\u05de\u05d5\u05e4\u05d7\u05ea\u05ea \u05d4\u05d5\u05e4\u05db\u05ea \u05d0\u05ea KANs \u05dc\u05de\u05ea\u05d0\u05d9\u05de\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d4\u05d3\u05d5\u05e8\u05e9\u05d5\u05ea \u05e9\u05e2\u05e8\u05d5\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05d1\u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd, \u05db\u05d2\u05d5\u05df \u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05d3\u05d9\u05e4\u05e8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea \u05d5\u05de\u05d9\u05d3\u05d5\u05dc \u05ea\u05d5\u05e4\u05e2\u05d5\u05ea \u05e4\u05d9\u05d6\u05d9\u05e7\u05dc\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea.\",\n  \" \u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d9\u05dd:\",\n  \"1. \u05de\u05d1\u05d7\u05e0\u05d9 \u05e8\u05d2\u05e8\u05e1\u05d9\u05d9\u05ea \u05ea\u05d3\u05e8\u05d9\u05dd: KANs \u05de\u05e6\u05dc\u05d9\u05d7\u05d5\u05ea \u05dc\u05d4\u05ea\u05d0\u05d9\u05dd \u05e8\u05db\u05d9\u05d1\u05d9 \u05d2\u05dc \u05d1\u05ea\u05d3\u05e8 \u05d2\u05d1\u05d5\u05d4 \u05d1\u05d5-\u05d6\u05de\u05e0\u05d9\u05ea, \u05d1\u05e2\u05d5\u05d3 \u05e9-MLPs \u05de\u05e6\u05d9\u05d2\u05d5\u05ea \u05e7\u05e9\u05d9\u05d9\u05dd \u05de\u05ea\u05de\u05e9\u05db\u05d9\u05dd \u05e2\u05dd \u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05d2\u05dd \u05dc\u05d0\u05d7\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05de\u05d5\u05e9\u05da.\",\n  \"2. \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9 \u05e9\u05d3\u05d4 \u05d2\u05d0\u05d5\u05e1\u05d9 \u05d0\u05e7\u05e8\u05d0\u05d9: KANs \u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d4\u05df \u05e2\u05dc MLPs \u05d1\u05e7\u05d9\u05e8\u05d5\u05d1 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e9\u05e0\u05d3\u05d2\u05de\u05d5 \u05de\u05e9\u05d3\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05d9\u05dd \u05d2\u05e1\u05d9\u05dd, \u05de\u05d4 \u05e9\u05de\u05e2\u05d9\u05d3 \u05e2\u05dc \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05e1\u05ea\u05d2\u05dc\u05d5\u05ea \u05e2\u05d3\u05d9\u05e4\u05d4 \u05dc\u05de\u05d1\u05e0\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd.\",\n  \"3. \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea PDE: \u05d1\u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05e4\u05d5\u05d0\u05e1\u05d5\u05df \u05d1\u05ea\u05d3\u05e8 \u05d2\u05d1\u05d5\u05d4, KANs \u05de\u05e9\u05d9\u05d2\u05d5\u05ea \u05e9\u05d2\u05d9\u05d0\u05d5\u05ea \u05e0\u05de\u05d5\u05db\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d0\u05d5\u05e4\u05df \u05e2\u05e7\u05d1\u05d9 \u05d1\u05d4\u05e9\u05d5\u05d5\u05d0\u05d4 \u05dc-MLPs, \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05d4 \u05e2\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d9\u05e6\u05d9\u05d1\u05d9\u05dd \u05d2\u05dd \u05db\u05d0\u05e9\u05e8 \u05ea\u05d3\u05e8 \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05e2\u05d5\u05dc\u05d4.\",\n  \"\u05d8\u05db\u05e0\u05d9\u05e7\u05ea \u05d4\u05e8\u05d7\u05d1\u05ea \u05d2\u05e8\u05d9\u05d3(\u05e9\u05dc \u05d4\u05e1\u05e4\u05dc\u05d9\u05d9\u05df):\",\n  \"\u05d7\u05d9\u05d3\u05d5\u05e9 \u05d8\u05db\u05e0\u05d9 \u05d1\u05d5\u05dc\u05d8 \u05d4\u05e0\u05d3\u05d5\u05df \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05d5\u05d0 \u05d8\u05db\u05e0\u05d9\u05e7\u05ea \u05d4\u05e8\u05d7\u05d1\u05ea \u05d2\u05e8\u05d9\u05d3 \u05d4\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05ea \u05dc- KANs. \u05e9\u05d9\u05d8\u05d4 \u05d6\u05d5 \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05e2\u05d9\u05d3\u05d5\u05df \u05d4\u05d3\u05e8\u05d2\u05ea\u05d9 \u05e9\u05dc \u05d2\u05e8\u05d9\u05d3 \u05e9\u05dc \u05d4-spline \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df, \u05d4\u05de\u05d0\u05e4\u05e9\u05e8 \u05ea\u05d4\u05dc\u05d9\u05da \u05dc\u05de\u05d9\u05d3\u05d4 \u05d9\u05e2\u05d9\u05dc \u05d9\u05d5\u05ea\u05e8. \u05d2\u05d9\u05e9\u05ea \u05d4\u05e8\u05d7\u05d1\u05ea \u05d4\u05d2\u05e8\u05d9\u05d3 \u05de\u05e4\u05d7\u05d9\u05ea\u05d4 \u05d0\u05ea \u05d4\u05e1\u05d9\u05db\u05d5\u05e0\u05d9\u05dd \u05dc-overfitting \u05d5\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05d4\u05db\u05dc\u05dc\u05d4 \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea, \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05db\u05d0\u05e9\u05e8 \u05de\u05ea\u05de\u05d5\u05d3\u05d3\u05d9\u05dd \u05e2\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d0\u05d5 \u05de\u05e2\u05e8\u05db\u05d9 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05d1\u05e2\u05dc\u05d9 \u05d3\u05d2\u05d9\u05de\u05d4 \u05d7\u05e1\u05e8\u05d4.\",\n  \"\u05e1\u05d9\u05db\u05d5\u05dd:\",\n  \"\u05e2\u05d1\u05d5\u05d3\u05d4 \u05d6\u05d5 \u05de\u05d1\u05e1\u05e1\u05ea \u05d0\u05ea KANs \u05db\u05d7\u05dc\u05d5\u05e4\u05d4 \u05d7\u05d6\u05e7\u05d4 \u05d5\u05d9\u05e2\u05d9\u05dc\u05d4 \u05dc\u05e8\u05e9\u05ea\u05d5\u05ea MLPs, \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d1\u05d7\u05d9\u05e9\u05d5\u05d1 \u05de\u05d3\u05e2\u05d9. \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05ea\u05de\u05d5\u05d3\u05d3\u05d5\u05ea \u05e2\u05dd \u05d4\u05d8\u05d9\u05d4 \u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea, \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1, \u05d5\u05e0\u05d9\u05e6\u05d5\u05dc \u05e9\u05d9\u05d8\u05d5\u05ea \u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05d3\u05e4\u05d8\u05d9\u05d1\u05d9\u05d5\u05ea, \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e1\u05e4\u05e7\u05d9\u05dd \u05e8\u05d0\u05d9\u05d5\u05ea \u05de\u05e9\u05db\u05e0\u05e2\u05d5\u05ea \u05dc\u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc \u05e9\u05dc KANs \u05dc\u05e2\u05dc\u05d5\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d4\u05df \u05e2\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05de\u05e1\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea \u05d1\u05d9\u05d9\u05e9\u05d5\u05de\u05d9\u05dd \u05d4\u05d3\u05d5\u05e8\u05e9\u05d9\u05dd \u05dc\u05de\u05d9\u05d3\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d5\u05de\u05e6\u05d9\u05d2\u05d5\u05ea \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1 \u05de\u05e9\u05d5\u05e4\u05e8\u05d5\u05ea. \u05d4\u05de\u05e1\u05d2\u05e8\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d1\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e2\u05dd \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9\u05dd \u05de\u05e7\u05d9\u05e4\u05d9\u05dd \u05d4\u05d5\u05e4\u05db\u05ea \u05de\u05d0\u05de\u05e8 \u05d6\u05d4 \u05dc\u05ea\u05e8\u05d5\u05de\u05d4 \u05d7\u05e9\u05d5\u05d1\u05d4 \u05dc\u05de\u05d7\u05e7\u05e8 \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd.\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nif (count !== 29) {\n  throw new Error(\"Unexpected paragraph count: \" + count);\n}\n\n// 1) First paragraph: title line + line break + paper title.\nparagraphs.items[0].insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 10.02.25\\u000bOn the expressiveness and spectral bias of KANs\", \"Replace\");\n\n// 2) Paragraphs 2-17 (zero-based indices 1-16): replace text in place.\nfor (let i = 0; i < middleTexts.length; i++) {\n  paragraphs.items[1 + i].insertText(middleTexts[i], \"Replace\");\n}\n\n// 3) Delete the old paragraphs that followed (zero-based indices 17-27,\n//    11 paragraphs) - delete from the end backwards so earlier indices\n//    stay valid.\nfor (let i = 27; i >= 17; i--) {\n  paragraphs.items[i].delete();\n}\n\n// 4) Last paragraph (the arxiv link): replace with the new URL.\nparagraphs.items[17].insertText(\"https://arxiv.org/abs/2410.01803\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is pre-seeded as $d below.\n#\n# Summary of the change (per the diff):\n#  - Paragraph 1 (title block): date 12.02.25 -> 10.02.25, and the paper\n#    title line (after the manual line break) is replaced with the new title.\n#  - Paragraphs 2-17 (the body of the review) get their text replaced in\n#    place, one-for-one, with the new \"KANs\" review content.\n#  - Paragraphs 18-28 (old content about State Collapse / SC strategies /\n#    summary) are removed entirely (11 paragraphs deleted).\n#  - The final paragraph (the arxiv link) keeps its position but its text\n#    changes to the new URL.\n\n$d = $word.ActiveDocument\n\nif ($d.Paragraphs.Count -ne 29) {\n    throw \"Unexpected paragraph count: $($d.Paragraphs.Count)\"\n}\n\n# 1) First paragraph: title line + manual line break (`v) + paper title.\n$d.Paragraphs(1).Range.Text = '\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 10.02.25' + \"`v\" + 'On the expressiveness and spectral bias of KANs'\n\n# 2) Paragraphs 2-17: replace text in place, one-for-one.\n$middleTexts = @(\n    '\u05de\u05d1\u05d5\u05d0:',\n    '\u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05d0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05de\u05e6\u05d9\u05d2 \u05d7\u05e7\u05d9\u05e8\u05d4 \u05de\u05e2\u05de\u05d9\u05e7\u05d4 \u05e9\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3 (KANs), \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d7\u05d3\u05e9\u05e0\u05d9\u05ea \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc \u05de\u05e9\u05e4\u05d8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e9\u05d5\u05d5\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05de\u05d3\u05d5\u05e7\u05d3\u05e7 \u05d1\u05d9\u05df KANs \u05dc\u05d1\u05d9\u05df \u05e8\u05e9\u05ea\u05d5\u05ea MLPs \u05de\u05e1\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea, \u05d4\u05df \u05de\u05d1\u05d7\u05d9\u05e0\u05d4 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d5\u05d4\u05df \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea, \u05ea\u05d5\u05da \u05d4\u05ea\u05de\u05e7\u05d3\u05d5\u05ea \u05d1\u05d4\u05d9\u05d1\u05d8\u05d9\u05dd \u05db\u05de\u05d5 \u05d0\u05e7\u05e1\u05e4\u05e8\u05e1\u05d9\u05d1\u05e0\u05e1, \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d5\u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05ea \u05d0\u05d9\u05de\u05d5\u05df. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d1\u05e1\u05e1 \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05e8\u05db\u05d6\u05d9\u05d5\u05ea \u05d5\u05de\u05d0\u05de\u05ea \u05d0\u05d5\u05ea\u05df \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9\u05dd, \u05d5\u05d1\u05db\u05da \u05de\u05d4\u05d5\u05d5\u05d4 \u05ea\u05e8\u05d5\u05de\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05dc\u05ea\u05db\u05e0\u05d5\u05df \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1 \u05e9\u05d5\u05e0\u05d5\u05ea.',\n    '\u05d0\u05e7\u05e1\u05e4\u05e8\u05e1\u05d9\u05d1\u05e0\u05e1:',\n    '\u05d4\u05d9\u05e9\u05d2 \u05de\u05e8\u05db\u05d6\u05d9 \u05e9\u05dc \u05e2\u05d1\u05d5\u05d3\u05d4 \u05d6\u05d5 \u05d4\u05d5\u05d0 \u05d4\u05d4\u05d5\u05db\u05d7\u05d4 \u05d4\u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e9- KANs \u05d4\u05df \u05d1\u05e2\u05dc\u05d5\u05ea \u05d0\u05e7\u05e1\u05e4\u05e8\u05e1\u05d9\u05d1\u05e0\u05e1 \u05dc\u05e4\u05d7\u05d5\u05ea \u05db\u05de\u05d5 MLPs. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05db\u05dc MLP \u05de\u05d1\u05d5\u05e1\u05e1\u05ea ReLU \u05e0\u05d9\u05ea\u05df \u05f4\u05dc\u05de\u05e4\u05d5\u05ea\u05f4 \u05dc\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea KAN \u05de\u05e7\u05d1\u05d9\u05dc\u05d4, \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05d4 \u05e2\u05dc \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d5\u05dc\u05dc\u05d0 \u05d4\u05d2\u05d3\u05dc\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05d1\u05d2\u05d5\u05d3\u05dc \u05d4\u05e8\u05e9\u05ea. \u05de\u05e0\u05d2\u05d3, \u05d1\u05e2\u05d5\u05d3 \u05e9-KANs \u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d2\u05dd \u05e2\u05dc \u05d9\u05d3\u05d9 MLPs, \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e6\u05d9\u05d4 \u05d6\u05d5 \u05db\u05e8\u05d5\u05db\u05d4 \u05d1\u05e2\u05dc\u05d5\u05ea \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea: \u05de\u05e1\u05e4\u05e8 \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d2\u05d3\u05dc \u05e2\u05dd \u05d2\u05d5\u05d3\u05dc \u05d2\u05e8\u05d9\u05d3 (\u05de\u05e1\u05e4\u05e8 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05e2\u05d5\u05d2\u05df \u05d1\u05e1\u05e4\u05dc\u05d9\u05d9\u05df) \u05e9\u05dc \u05d4-KAN. \u05de\u05de\u05e6\u05d0 \u05d6\u05d4 \u05de\u05e8\u05de\u05d6 \u05e9-KANs \u05e2\u05e9\u05d5\u05d9\u05d5\u05ea \u05dc\u05d4\u05e6\u05d9\u05e2 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d9\u05e2\u05d9\u05dc\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05e2\u05d1\u05d5\u05e8 \u05e1\u05d5\u05d2\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea, \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05db\u05d0\u05e9\u05e8 \u05e0\u05e2\u05e9\u05d4 \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05de\u05d1\u05e0\u05d9 \u05d2\u05e8\u05d9\u05d3 \u05e2\u05d3\u05d9\u05e0\u05d9\u05dd.',\n    '\u05d4\u05de\u05d7\u05e7\u05e8 \u05de\u05e0\u05e6\u05dc \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05e7\u05d9\u05d9\u05de\u05d5\u05ea \u05e2\u05d1\u05d5\u05e8 MLPs \u05db\u05d3\u05d9 \u05dc\u05e7\u05d1\u05d5\u05e2 \u05e7\u05e6\u05d1\u05d9 \u05e7\u05d9\u05e8\u05d5\u05d1 \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05d1\u05d5\u05e8 KANs \u05d1\u05de\u05e8\u05d7\u05d1\u05d9\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05d0\u05d5\u05e0\u05dc\u05d9\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05db\u05de\u05d5 \u05de\u05e8\u05d7\u05d1 \u05e1\u05d5\u05d1\u05d5\u05dc\u05d1. \u05d4\u05d5\u05d0 \u05de\u05d3\u05d2\u05d9\u05dd \u05e9-KANs \u05de\u05e9\u05d9\u05d2\u05d5\u05ea \u05e7\u05e6\u05d1\u05d9 \u05e7\u05d9\u05e8\u05d5\u05d1 \u05d3\u05d5\u05de\u05d9\u05dd \u05d0\u05d5 \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d0\u05e9\u05e8 MLPs \u05d1\u05e9\u05e2\u05e8\u05d5\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea, \u05de\u05d4 \u05e9\u05de\u05d7\u05d6\u05e7 \u05d0\u05ea \u05d7\u05d5\u05e1\u05e0\u05df \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9.',\n    '\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05d8\u05d9\u05d9\u05ea \u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea (spectral bias):',\n    '\u05d0\u05d7\u05d3 \u05d4\u05d4\u05d1\u05d3\u05dc\u05d9\u05dd \u05d4\u05de\u05e8\u05db\u05d6\u05d9\u05d9\u05dd \u05d1\u05d9\u05df KANs \u05dc-MLPs \u05d4\u05de\u05d5\u05d3\u05d2\u05e9\u05d9\u05dd \u05d1\u05de\u05d0\u05de\u05e8 \u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d4\u05d4\u05d1\u05d3\u05dc \u05d1\u05d4\u05d8\u05d9\u05d4 \u05d4\u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea \u05e9\u05dc\u05d4\u05dd - \u05ea\u05d5\u05e4\u05e2\u05d4 \u05e9\u05d1\u05d4 \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05e0\u05d5\u05d8\u05d5\u05ea \u05dc\u05dc\u05de\u05d5\u05d3 \u05ea\u05d7\u05d9\u05dc\u05d4 \u05d1\u05ea\u05d3\u05e8\u05d9\u05dd \u05e0\u05de\u05d5\u05db\u05d9\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05d2\u05d9\u05dd \u05e0\u05d9\u05ea\u05d5\u05d7 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05d5\u05d0\u05de\u05e4\u05d9\u05e8\u05d9 \u05de\u05e4\u05d5\u05e8\u05d8, \u05d4\u05de\u05e8\u05d0\u05d4 \u05e9- KANs \u05e1\u05d5\u05d1\u05dc\u05d5\u05ea \u05e4\u05d7\u05d5\u05ea \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05de\u05d4\u05d8\u05d9\u05d4 \u05d6\u05d5.',\n    '\u05d4\u05d1\u05d3\u05dc \u05d6\u05d4 \u05de\u05d9\u05d5\u05d7\u05e1 \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05d4-B-spline \u05d5\u05dc\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d4\u05e7\u05d5\u05de\u05e4\u05d5\u05d6\u05d9\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9\u05ea \u05e9\u05dc KANs, \u05d4\u05de\u05d0\u05e4\u05e9\u05e8\u05d5\u05ea \u05dc\u05d4\u05df \u05dc\u05dc\u05de\u05d5\u05d3 \u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4 \u05d1\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05e8\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05e6\u05d9\u05e2\u05d5\u05ea \u05e9\u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05ea \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc KANs \u05e8\u05d3\u05d5\u05d3\u05d5\u05ea \u05d0\u05d7\u05d9\u05d3\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d9\u05d7\u05e1 \u05dc\u05ea\u05d3\u05e8\u05d9\u05dd \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05d1\u05d4\u05e9\u05d5\u05d5\u05d0\u05d4 \u05dc-MLPs, \u05e9\u05d1\u05d4\u05df \u05e0\u05e6\u05e4\u05d9\u05ea \u05d4\u05ea\u05db\u05e0\u05e1\u05d5\u05ea \u05de\u05d4\u05d9\u05e8\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05ea\u05d3\u05e8\u05d9\u05dd \u05e0\u05de\u05d5\u05db\u05d9\u05dd. \u05d4\u05d4\u05d8\u05d9\u05d4 \u05d4\u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea \u05d4\u05de\u05d5\u05e4\u05d7\u05ea\u05ea \u05d4\u05d5\u05e4\u05db\u05ea \u05d0\u05ea KANs \u05dc\u05de\u05ea\u05d0\u05d9\u05de\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d4\u05d3\u05d5\u05e8\u05e9\u05d5\u05ea \u05e9\u05e2\u05e8\u05d5\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05d1\u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd, \u05db\u05d2\u05d5\u05df \u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05d3\u05d9\u05e4\u05e8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea \u05d5\u05de\u05d9\u05d3\u05d5\u05dc \u05ea\u05d5\u05e4\u05e2\u05d5\u05ea \u05e4\u05d9\u05d6\u05d9\u05e7\u05dc\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea.',\n    ' \u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d9\u05dd:',\n    '1. \u05de\u05d1\u05d7\u05e0\u05d9 \u05e8\u05d2\u05e8\u05e1\u05d9\u05d9\u05ea \u05ea\u05d3\u05e8\u05d9\u05dd: KANs \u05de\u05e6\u05dc\u05d9\u05d7\u05d5\u05ea \u05dc\u05d4\u05ea\u05d0\u05d9\u05dd \u05e8\u05db\u05d9\u05d1\u05d9 \u05d2\u05dc \u05d1\u05ea\u05d3\u05e8 \u05d2\u05d1\u05d5\u05d4 \u05d1\u05d5-\u05d6\u05de\u05e0\u05d9\u05ea, \u05d1\u05e2\u05d5\u05d3 \u05e9-MLPs \u05de\u05e6\u05d9\u05d2\u05d5\u05ea \u05e7\u05e9\u05d9\u05d9\u05dd \u05de\u05ea\u05de\u05e9\u05db\u05d9\u05dd \u05e2\u05dd \u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05d2\u05dd \u05dc\u05d0\u05d7\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05de\u05d5\u05e9\u05da.',\n    '2. \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9 \u05e9\u05d3\u05d4 \u05d2\u05d0\u05d5\u05e1\u05d9 \u05d0\u05e7\u05e8\u05d0\u05d9: KANs \u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d4\u05df \u05e2\u05dc MLPs \u05d1\u05e7\u05d9\u05e8\u05d5\u05d1 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e9\u05e0\u05d3\u05d2\u05de\u05d5 \u05de\u05e9\u05d3\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05d9\u05dd \u05d2\u05e1\u05d9\u05dd, \u05de\u05d4 \u05e9\u05de\u05e2\u05d9\u05d3 \u05e2\u05dc \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05e1\u05ea\u05d2\u05dc\u05d5\u05ea \u05e2\u05d3\u05d9\u05e4\u05d4 \u05dc\u05de\u05d1\u05e0\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd.',\n    '3. \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea PDE: \u05d1\u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05e4\u05d5\u05d0\u05e1\u05d5\u05df \u05d1\u05ea\u05d3\u05e8 \u05d2\u05d1\u05d5\u05d4, KANs \u05de\u05e9\u05d9\u05d2\u05d5\u05ea \u05e9\u05d2\u05d9\u05d0\u05d5\u05ea \u05e0\u05de\u05d5\u05db\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d0\u05d5\u05e4\u05df \u05e2\u05e7\u05d1\u05d9 \u05d1\u05d4\u05e9\u05d5\u05d5\u05d0\u05d4 \u05dc-MLPs, \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05d4 \u05e2\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d9\u05e6\u05d9\u05d1\u05d9\u05dd \u05d2\u05dd \u05db\u05d0\u05e9\u05e8 \u05ea\u05d3\u05e8 \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05e2\u05d5\u05dc\u05d4.',\n    '\u05d8\u05db\u05e0\u05d9\u05e7\u05ea \u05d4\u05e8\u05d7\u05d1\u05ea \u05d2\u05e8\u05d9\u05d3(\u05e9\u05dc \u05d4\u05e1\u05e4\u05dc\u05d9\u05d9\u05df):',\n    '\u05d7\u05d9\u05d3\u05d5\u05e9 \u05d8\u05db\u05e0\u05d9 \u05d1\u05d5\u05dc\u05d8 \u05d4\u05e0\u05d3\u05d5\u05df \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05d5\u05d0 \u05d8\u05db\u05e0\u05d9\u05e7\u05ea \u05d4\u05e8\u05d7\u05d1\u05ea \u05d2\u05e8\u05d9\u05d3 \u05d4\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05ea \u05dc- KANs. \u05e9\u05d9\u05d8\u05d4 \u05d6\u05d5 \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05e2\u05d9\u05d3\u05d5\u05df \u05d4\u05d3\u05e8\u05d2\u05ea\u05d9 \u05e9\u05dc \u05d2\u05e8\u05d9\u05d3 \u05e9\u05dc \u05d4-spline \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df, \u05d4\u05de\u05d0\u05e4\u05e9\u05e8 \u05ea\u05d4\u05dc\u05d9\u05da \u05dc\u05de\u05d9\u05d3\u05d4 \u05d9\u05e2\u05d9\u05dc \u05d9\u05d5\u05ea\u05e8. \u05d2\u05d9\u05e9\u05ea \u05d4\u05e8\u05d7\u05d1\u05ea \u05d4\u05d2\u05e8\u05d9\u05d3 \u05de\u05e4\u05d7\u05d9\u05ea\u05d4 \u05d0\u05ea \u05d4\u05e1\u05d9\u05db\u05d5\u05e0\u05d9\u05dd \u05dc-overfitting \u05d5\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05d4\u05db\u05dc\u05dc\u05d4 \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea, \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05db\u05d0\u05e9\u05e8 \u05de\u05ea\u05de\u05d5\u05d3\u05d3\u05d9\u05dd \u05e2\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d0\u05d5 \u05de\u05e2\u05e8\u05db\u05d9 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05d1\u05e2\u05dc\u05d9 \u05d3\u05d2\u05d9\u05de\u05d4 \u05d7\u05e1\u05e8\u05d4.',\n    '\u05e1\u05d9\u05db\u05d5\u05dd:',\n    '\u05e2\u05d1\u05d5\u05d3\u05d4 \u05d6\u05d5 \u05de\u05d1\u05e1\u05e1\u05ea \u05d0\u05ea KANs \u05db\u05d7\u05dc\u05d5\u05e4\u05d4 \u05d7\u05d6\u05e7\u05d4 \u05d5\u05d9\u05e2\u05d9\u05dc\u05d4 \u05dc\u05e8\u05e9\u05ea\u05d5\u05ea MLPs, \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d1\u05d7\u05d9\u05e9\u05d5\u05d1 \u05de\u05d3\u05e2\u05d9. \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05ea\u05de\u05d5\u05d3\u05d3\u05d5\u05ea \u05e2\u05dd \u05d4\u05d8\u05d9\u05d4 \u05e1\u05e4\u05e7\u05d8\u05e8\u05dc\u05d9\u05ea, \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1, \u05d5\u05e0\u05d9\u05e6\u05d5\u05dc \u05e9\u05d9\u05d8\u05d5\u05ea \u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05d3\u05e4\u05d8\u05d9\u05d1\u05d9\u05d5\u05ea, \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e1\u05e4\u05e7\u05d9\u05dd \u05e8\u05d0\u05d9\u05d5\u05ea \u05de\u05e9\u05db\u05e0\u05e2\u05d5\u05ea \u05dc\u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc \u05e9\u05dc KANs \u05dc\u05e2\u05dc\u05d5\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d4\u05df \u05e2\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05de\u05e1\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea \u05d1\u05d9\u05d9\u05e9\u05d5\u05de\u05d9\u05dd \u05d4\u05d3\u05d5\u05e8\u05e9\u05d9\u05dd \u05dc\u05de\u05d9\u05d3\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05ea\u05d3\u05e8\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d5\u05de\u05e6\u05d9\u05d2\u05d5\u05ea \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1 \u05de\u05e9\u05d5\u05e4\u05e8\u05d5\u05ea. \u05d4\u05de\u05e1\u05d2\u05e8\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d1\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e2\u05dd \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9\u05dd \u05de\u05e7\u05d9\u05e4\u05d9\u05dd \u05d4\u05d5\u05e4\u05db\u05ea \u05de\u05d0\u05de\u05e8 \u05d6\u05d4 \u05dc\u05ea\u05e8\u05d5\u05de\u05d4 \u05d7\u05e9\u05d5\u05d1\u05d4 \u05dc\u05de\u05d7\u05e7\u05e8 \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd.'\n)\nfor ($i = 0; $i -lt $middleTexts.Count; $i++) {\n    $d.Paragraphs(2 + $i).Range.Text = $middleTexts[$i]\n}\n\n# 3) Delete the old paragraphs that followed (1-based paragraphs 18-28,\n#    11 paragraphs) via a single Range delete.\n$startPara = $d.Paragraphs(18)\n$endPara = $d.Paragraphs(28)\n$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$rng.Delete()\n\n# 4) Last paragraph (the arxiv link): replace with the new URL.\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = 'https://arxiv.org/abs/2410.01803'\n\n"}
